# Applies the cryptos.xlsx data refresh described in the commit message
# "Updated cryptos list on Wed Aug 23 03:56:02 UTC 2023 with GitHub Actions".
# Every price/volume cell in the sheet is stored as literal text (not a
# real number), including values that happen to look numeric (e.g. "20.64")
# -- so every write below pins the cell to the Text number format before
# setting .Value (stops Excel's COM layer from silently re-typing strings
# like "0.5230"/"64.80" as numbers and dropping trailing zeros), then
# resets the style back to "Normal" so no stray numFmt/style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 "26.047.32"
Set-TextCell 2 5 "  -0.27%  "

# Row 3
Set-TextCell 3 4 "1.634.75"
Set-TextCell 3 5 "  -1.95%  "

# Row 4
Set-TextCell 4 5 "  +0.00%  "

# Row 5
Set-TextCell 5 4 "212.84"
Set-TextCell 5 5 "  +1.02%  "

# Row 6
Set-TextCell 6 4 "0.5230"
Set-TextCell 6 5 "  -0.34%  "

# Row 7
Set-TextCell 7 5 "  +0.00%  "

# Row 8
Set-TextCell 8 4 "0.2592"
Set-TextCell 8 5 "  -1.33%  "

# Row 9
Set-TextCell 9 4 "0.06281"
Set-TextCell 9 5 "  -0.10%  "

# Row 10
Set-TextCell 10 4 "20.64"
Set-TextCell 10 5 "  -2.37%  "

# Row 11
Set-TextCell 11 4 "0.07653"
Set-TextCell 11 5 "  +1.56%  "

# Row 12
Set-TextCell 12 4 "1.633.09"
Set-TextCell 12 5 "  -2.02%  "

# Row 13
Set-TextCell 13 4 "4.406"
Set-TextCell 13 5 "  -0.70%  "

# Row 14
Set-TextCell 14 4 "1.858.82"
Set-TextCell 14 5 "  -1.96%  "

# Row 15
Set-TextCell 15 4 "0.5505"
Set-TextCell 15 5 "  -0.74%  "

# Row 16
Set-TextCell 16 4 "0.0₅8165"
Set-TextCell 16 5 "  +2.88%  "

# Row 17
Set-TextCell 17 4 "64.80"
Set-TextCell 17 5 "  -2.85%  "

# Row 18
Set-TextCell 18 4 "26.045.82"
Set-TextCell 18 5 "  -0.34%  "

# Row 19
Set-TextCell 19 5 "  +0.01%  "

# Row 20
Set-TextCell 20 4 "4.677"
Set-TextCell 20 5 "  -1.10%  "

# Row 21
Set-TextCell 21 4 "187.81"
Set-TextCell 21 5 "  +0.77%  "

# Row 22
Set-TextCell 22 4 "10.15"
Set-TextCell 22 5 "  -1.76%  "

# Row 23
Set-TextCell 23 4 "6.135"

# Row 24
Set-TextCell 24 5 "  +0.10%  "

# Row 25
Set-TextCell 25 4 "145.05"
Set-TextCell 25 5 "  -3.01%  "

# Row 26
Set-TextCell 26 5 "  -2.78%  "

# Row 27
Set-TextCell 27 4 "7.386"
Set-TextCell 27 5 "  -1.21%  "

# Row 28
Set-TextCell 28 5 "  -1.04%  "

# Row 29
Set-TextCell 29 4 "1.397"
Set-TextCell 29 5 "  +3.15%  "

# Row 30
Set-TextCell 30 4 "0.05933"
Set-TextCell 30 5 "  -5.16%  "

# Row 31
Set-TextCell 31 4 "1.253"
Set-TextCell 31 5 "  -2.00%  "

# Row 32
Set-TextCell 32 5 "  -1.98%  "

# Row 33
Set-TextCell 33 4 "3.403"
Set-TextCell 33 5 "  -0.22%  "

# Row 34
Set-TextCell 34 4 "1.634"
Set-TextCell 34 5 "  +0.18%  "

# Row 35
Set-TextCell 35 4 "0.9817"
Set-TextCell 35 5 "  -1.66%  "

# Row 36
Set-TextCell 36 5 "  -0.69%  "

# Row 37
Set-TextCell 37 4 "2.761"
Set-TextCell 37 5 "  +1.18%  "

# Row 38
Set-TextCell 38 4 "0.5698"
Set-TextCell 38 5 "  -5.38%  "

# Row 39
Set-TextCell 39 4 "0.01613"
Set-TextCell 39 5 "  -0.14%  "

# Row 40
Set-TextCell 40 4 "0.8513"
Set-TextCell 40 5 "  -2.13%  "

# Row 41
Set-TextCell 41 2 "FraxShare"
Set-TextCell 41 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell 41 4 "5.739"
Set-TextCell 41 5 "  -6.11%  "

# Row 42
Set-TextCell 42 2 "PaxDollar"
Set-TextCell 42 3 "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell 42 4 "1.001"
Set-TextCell 42 5 "  -0.14%  "

# Row 43
Set-TextCell 43 4 "1.031.90"
Set-TextCell 43 5 "  -6.64%  "

# Row 44
Set-TextCell 44 4 "100.30"
Set-TextCell 44 5 "  +0.30%  "

# Row 45
Set-TextCell 45 4 "1.785.23"
Set-TextCell 45 5 "  -1.86%  "

# Row 46
Set-TextCell 46 2 "BabyDogeCoin"
Set-TextCell 46 3 "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell 46 4 "0.0₈110"
Set-TextCell 46 5 "  -1.53%  "

# Row 47
Set-TextCell 47 2 "Aave"
Set-TextCell 47 3 "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell 47 4 "55.68"
Set-TextCell 47 5 "  +0.57%  "

# Row 48
Set-TextCell 48 4 "1.001"
Set-TextCell 48 5 "  -0.03%  "

# Row 49
Set-TextCell 49 4 "8.018"
Set-TextCell 49 5 "  -0.40%  "

# Row 50
Set-TextCell 50 4 "0.05164"
Set-TextCell 50 5 "  -1.37%  "

# Row 51
Set-TextCell 51 5 "  -0.56%  "
